# chore: update Sheets via scheduled runner
# Refresh Universalis market-price snapshots (currentAveragePrice* /
# LevePrice* / LeveProfit* columns, H:N) across the gathering/crafting
# job sheets. Values below are the latest scrape for the affected leves.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1470.3334
$ws.Range("I113").Value = 1949.5
$ws.Range("J113").Value = 512
$ws.Range("K113").Value = 1949.5
$ws.Range("L113").Value = 512
$ws.Range("M113").Value = 1304.5
$ws.Range("N113").Value = -7020
$ws.Range("H116").Value = 9998
$ws.Range("J116").Value = 9998
$ws.Range("L116").Value = 9998
$ws.Range("N116").Value = -16882
$ws.Range("H137").Value = 1120221.8
$ws.Range("J137").Value = 19666.666
$ws.Range("L137").Value = 58999.99800000001
$ws.Range("N137").Value = -64099.99800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 35000
$ws.Range("J96").Value = 35000
$ws.Range("L96").Value = 35000
$ws.Range("N96").Value = -40492
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("H122").Value = 2228.1428
$ws.Range("I122").Value = 1649.5
$ws.Range("J122").Value = 2999.6667
$ws.Range("K122").Value = 4948.5
$ws.Range("L122").Value = 8999.000100000001
$ws.Range("M122").Value = -2498.5
$ws.Range("N122").Value = -13899.0001
$ws.Range("H123").Value = 20000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 20000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -29800
$ws.Range("H124").Value = 58623.875
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 58623.875
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 58623.875
$ws.Range("N124").Value = -68443.875
$ws.Range("H125").Value = 183332.33
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 183332.33
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 183332.33
$ws.Range("N125").Value = -193172.33
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("H127").Value = 18889
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 18889
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 18889
$ws.Range("N127").Value = -28809
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H130").Value = 70000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 70000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 70000
$ws.Range("N130").Value = -80040
$ws.Range("H131").Value = 50000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 50000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H132").Value = 13331.667
$ws.Range("I132").Value = 9997.5
$ws.Range("J132").Value = 20000
$ws.Range("K132").Value = 29992.5
$ws.Range("L132").Value = 60000
$ws.Range("M132").Value = -27462.5
$ws.Range("N132").Value = -65060
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 34497.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 34497.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 34497.5
$ws.Range("N135").Value = -44637.5
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 50000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 50000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280
$ws.Range("H139").Value = 71000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 71000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 71000
$ws.Range("N139").Value = -81280
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1187.5
$ws.Range("I16").Value = 1187.5
$ws.Range("K16").Value = 1187.5
$ws.Range("M16").Value = -900.5
$ws.Range("H105").Value = 1399.75
$ws.Range("I105").Value = 1399.75
$ws.Range("K105").Value = 1399.75
$ws.Range("M105").Value = 347.25
$ws.Range("H113").Value = 1187.5
$ws.Range("I113").Value = 1187.5
$ws.Range("K113").Value = 1187.5
$ws.Range("M113").Value = 982.5
$ws.Range("H122").Value = 1637.25
$ws.Range("I122").Value = 1116.3334
$ws.Range("K122").Value = 3349.0002
$ws.Range("M122").Value = -899.0001999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 91.833336
$ws.Range("J23").Value = 100.2
$ws.Range("L23").Value = 300.6
$ws.Range("N23").Value = -770.6
$ws.Range("H33").Value = 46.5
$ws.Range("I33").Value = 67.5
$ws.Range("J33").Value = 36
$ws.Range("K33").Value = 405
$ws.Range("L33").Value = 216
$ws.Range("M33").Value = -122
$ws.Range("N33").Value = -782
$ws.Range("H80").Value = 14002
$ws.Range("I80").Value = 14002
$ws.Range("K80").Value = 42006
$ws.Range("M80").Value = -41070
$ws.Range("H83").Value = 14002
$ws.Range("I83").Value = 14002
$ws.Range("K83").Value = 126018
$ws.Range("M83").Value = -121338
$ws.Range("H107").Value = 779.8
$ws.Range("I107").Value = 600
$ws.Range("K107").Value = 1800
$ws.Range("L107").Value = 2699.0001
$ws.Range("M107").Value = 120
$ws.Range("N107").Value = -6539.0001
$ws.Range("H122").Value = 175.5
$ws.Range("I122").Value = 175.5
$ws.Range("K122").Value = 1579.5
$ws.Range("M122").Value = 870.5
$ws.Range("H136").Value = 5749.5
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 9999
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 29997
$ws.Range("M136").Value = 600
$ws.Range("N136").Value = -40197

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 720
$ws.Range("I22").Value = 670
$ws.Range("J22").Value = 820
$ws.Range("K22").Value = 670
$ws.Range("L22").Value = 820
$ws.Range("M22").Value = -375
$ws.Range("N22").Value = -1410
$ws.Range("H27").Value = 720
$ws.Range("I27").Value = 670
$ws.Range("J27").Value = 820
$ws.Range("K27").Value = 670
$ws.Range("L27").Value = 820
$ws.Range("M27").Value = -563
$ws.Range("N27").Value = -1034
$ws.Range("H35").Value = 682.5
$ws.Range("I35").Value = 682.5
$ws.Range("K35").Value = 682.5
$ws.Range("M35").Value = -346.5
$ws.Range("H100").Value = 1501.5
$ws.Range("I100").Value = 1501.5
$ws.Range("K100").Value = 1501.5
$ws.Range("M100").Value = -960.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2808
$ws.Range("I122").Value = 2369.6
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7108.799999999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4658.799999999999
$ws.Range("N122").Value = -19900

